# Add "Possible conductor material" values into column G of Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# G25 already exists (empty, with style 3 inherited from row 25's "heavy" styling).
# Re-format it to the plain style ("1") used throughout the column before giving
# it a value, by pasting the formats from a cell that already carries style 1.
$ws.Range("A2").Copy() | Out-Null
$ws.Range("G25").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Rows which get "Nb3Sn" in column G
$nb3snRows = @(2,3,4,5,7,8,9,11,12,14,15,17,18,19,21,22,23,25,26,27,29,30,31,34,35)
foreach ($r in $nb3snRows) {
    $ws.Range("G$r").Value = "Nb3Sn"
}

# Row 33 gets "REBCO" in column G
$ws.Range("G33").Value = "REBCO"

# Update two data values that changed
$ws.Range("C31").Value = 12.1
$ws.Range("C35").Value = 10.6

# Update the view state (selection) to match the saved file
$ws.Activate()
$ws.Range("C27").Select()
